$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 14374
$ws.Range("J32").Value = 16284.857
$ws.Range("L32").Value = 16284.857
$ws.Range("N32").Value = -16936.857

$ws.Range("H40").Value = 3630.4211
$ws.Range("I40").Value = 2936.0625
$ws.Range("J40").Value = 7333.6665
$ws.Range("K40").Value = 2936.0625
$ws.Range("L40").Value = 7333.6665
$ws.Range("M40").Value = -2761.0625
$ws.Range("N40").Value = -7683.6665

$ws.Range("H96").Value = 1248
$ws.Range("I96").Value = 887.44446
$ws.Range("K96").Value = 2662.33338
$ws.Range("M96").Value = -1289.33338

$ws.Range("H98").Value = 820.3077
$ws.Range("I98").Value = 820.3077
$ws.Range("K98").Value = 820.3077
$ws.Range("M98").Value = 677.6923

$ws.Range("H116").Value = 5318.8423
$ws.Range("I116").Value = 5091.647
$ws.Range("J116").Value = 7250
$ws.Range("K116").Value = 5091.647
$ws.Range("L116").Value = 7250
$ws.Range("M116").Value = -1649.647
$ws.Range("N116").Value = -14134

$ws.Range("H122").Value = 820.3077
$ws.Range("I122").Value = 820.3077
$ws.Range("K122").Value = 2460.9231
$ws.Range("M122").Value = -10.92309999999998

$ws.Range("H125").Value = 1829.7142
$ws.Range("I125").Value = 1747.4
$ws.Range("K125").Value = 15726.6
$ws.Range("M125").Value = -13266.6

$ws.Range("H137").Value = 3114.6
$ws.Range("I137").Value = 2828.4546
$ws.Range("J137").Value = 3901.5
$ws.Range("K137").Value = 8485.363799999999
$ws.Range("L137").Value = 11704.5
$ws.Range("M137").Value = -5935.363799999999
$ws.Range("N137").Value = -16804.5

$ws.Range("H138").Value = 2955.0588
$ws.Range("J138").Value = 3488.8235
$ws.Range("L138").Value = 10466.4705
$ws.Range("N138").Value = -20746.4705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3073.5264
$ws.Range("I105").Value = 2569.077
$ws.Range("K105").Value = 2569.077
$ws.Range("M105").Value = -822.0770000000002

$ws.Range("H134").Value = 31253136
$ws.Range("I134").Value = 33335346
$ws.Range("K134").Value = 100006038
$ws.Range("M134").Value = -100003503

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 728262
$ws.Range("I16").Value = 989811.9399999999
$ws.Range("K16").Value = 989811.9399999999
$ws.Range("M16").Value = -989524.9399999999

$ws.Range("H31").Value = 9211.281999999999
$ws.Range("I31").Value = 7037.8335
$ws.Range("J31").Value = 11582.318
$ws.Range("K31").Value = 7037.8335
$ws.Range("L31").Value = 11582.318
$ws.Range("M31").Value = -6742.8335
$ws.Range("N31").Value = -12172.318

$ws.Range("H34").Value = 9211.281999999999
$ws.Range("I34").Value = 7037.8335
$ws.Range("J34").Value = 11582.318
$ws.Range("K34").Value = 7037.8335
$ws.Range("L34").Value = 11582.318
$ws.Range("M34").Value = -6835.8335
$ws.Range("N34").Value = -11986.318

$ws.Range("H58").Value = 27785738
$ws.Range("I58").Value = 35723092
$ws.Range("J58").Value = 5006.5
$ws.Range("K58").Value = 35723092
$ws.Range("L58").Value = 5006.5
$ws.Range("M58").Value = -35722889
$ws.Range("N58").Value = -5412.5

$ws.Range("H113").Value = 728262
$ws.Range("I113").Value = 989811.9399999999
$ws.Range("K113").Value = 989811.9399999999
$ws.Range("M113").Value = -987641.9399999999

$ws.Range("H132").Value = 26317186
$ws.Range("I132").Value = 26317186
$ws.Range("K132").Value = 78951558
$ws.Range("M132").Value = -78949028

$ws.Range("H136").Value = 27785738
$ws.Range("I136").Value = 35723092
$ws.Range("J136").Value = 5006.5
$ws.Range("K136").Value = 107169276
$ws.Range("L136").Value = 15019.5
$ws.Range("M136").Value = -107166726
$ws.Range("N136").Value = -20119.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 4000

$ws.Range("H121").Value = 65325.05
$ws.Range("I121").Value = 127629
$ws.Range("J121").Value = 20013.092
$ws.Range("K121").Value = 382887
$ws.Range("L121").Value = 60039.276
$ws.Range("M121").Value = -381577
$ws.Range("N121").Value = -62659.276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7817194.5
$ws.Range("I132").Value = 9618462
$ws.Range("J132").Value = 11703.667
$ws.Range("K132").Value = 28855386
$ws.Range("L132").Value = 35111.001
$ws.Range("M132").Value = -28852856
$ws.Range("N132").Value = -40171.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2941.1428
$ws.Range("I40").Value = 2941.1428
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2941.1428
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2805.1428
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 923.8
$ws.Range("I46").Value = 846.9286
$ws.Range("K46").Value = 846.9286
$ws.Range("M46").Value = -658.9286

$ws.Range("H61").Value = 2906.7693
$ws.Range("I61").Value = 2639.1
$ws.Range("J61").Value = 3799
$ws.Range("K61").Value = 2639.1
$ws.Range("L61").Value = 3799
$ws.Range("M61").Value = -2437.1
$ws.Range("N61").Value = -4203

$ws.Range("H113").Value = 2906.7693
$ws.Range("I113").Value = 2639.1
$ws.Range("J113").Value = 3799
$ws.Range("K113").Value = 2639.1
$ws.Range("L113").Value = 3799
$ws.Range("M113").Value = -469.0999999999999
$ws.Range("N113").Value = -8139

$ws.Range("H122").Value = 4876.4814
$ws.Range("I122").Value = 5024.3477
$ws.Range("J122").Value = 4026.25
$ws.Range("K122").Value = 15073.0431
$ws.Range("L122").Value = 12078.75
$ws.Range("M122").Value = -12623.0431
$ws.Range("N122").Value = -16978.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H54").Value = 3070
$ws.Range("I54").Value = 3070
$ws.Range("K54").Value = 3070
$ws.Range("M54").Value = -2550

$ws.Range("H122").Value = 5256.2
$ws.Range("I122").Value = 6030.4165
$ws.Range("K122").Value = 18091.2495
$ws.Range("M122").Value = -15641.2495

$ws.Range("H132").Value = 19237278
$ws.Range("I132").Value = 31252586
$ws.Range("J132").Value = 12785.8
$ws.Range("K132").Value = 93757758
$ws.Range("L132").Value = 38357.39999999999
$ws.Range("M132").Value = -93755228
$ws.Range("N132").Value = -43417.39999999999
